$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column-D cells whose new value looks like a plain number need the
# cell format forced to Text first, otherwise Excel auto-converts the
# assigned string into a numeric value (dropping trailing zeros, losing
# the exact decimal representation, etc.). Values that already contain
# two "thousands" separators (e.g. "45.801.10") are never auto-parsed as
# numbers, so they do not need this treatment.
$textCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D17", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D37", "D39", "D40", "D42", "D44", "D45", "D46", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "45.801.10"
$ws.Range("E2").Value = "  -2.30%  "

$ws.Range("D3").Value = "2.368.55"
$ws.Range("E3").Value = "  +2.61%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "299.75"
$ws.Range("E5").Value = "  -1.52%  "

$ws.Range("D6").Value = "98.17"
$ws.Range("E6").Value = "  -3.66%  "

$ws.Range("D7").Value = "0.565"
$ws.Range("E7").Value = "  -1.18%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  -4.45%  "

$ws.Range("D10").Value = "34.17"
$ws.Range("E10").Value = "  -7.71%  "

$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  -2.18%  "

$ws.Range("D12").Value = "7.08"
$ws.Range("E12").Value = "  -4.74%  "

$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("D14").Value = "2.732.27"
$ws.Range("E14").Value = "  +2.73%  "

$ws.Range("D15").Value = "2.377.38"
$ws.Range("E15").Value = "  +2.87%  "

$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").Value = "13.67"
$ws.Range("E17").Value = "  -2.75%  "

$ws.Range("D18").Value = "45.747.13"
$ws.Range("E18").Value = "  -2.35%  "

$ws.Range("E19").Value = "  -8.48%  "

$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  -1.34%  "

$ws.Range("D22").Value = "66.67"
$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("D23").Value = "243.05"
$ws.Range("E23").Value = "  -2.61%  "

$ws.Range("D24").Value = "2.77"
$ws.Range("E24").Value = "  -5.78%  "

$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").Value = "1.92"
$ws.Range("E26").Value = "  -2.31%  "

$ws.Range("D27").Value = "38.87"
$ws.Range("E27").Value = "  -11.16%  "

$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -4.02%  "

$ws.Range("D29").Value = "9.69"
$ws.Range("E29").Value = "  -3.06%  "

$ws.Range("D30").Value = "20.88"
$ws.Range("E30").Value = "  +3.03%  "

$ws.Range("D31").Value = "3.71"
$ws.Range("E31").Value = "  +15.45%  "

$ws.Range("D32").Value = "2.75"
$ws.Range("E32").Value = "  -4.61%  "

$ws.Range("D33").Value = "5.53"
$ws.Range("E33").Value = "  -5.30%  "

$ws.Range("D34").Value = "147.16"
$ws.Range("E34").Value = "  -0.41%  "

$ws.Range("E35").Value = "  -5.08%  "

$ws.Range("E36").Value = "  -1.10%  "

$ws.Range("D37").Value = "1.91"
$ws.Range("E37").Value = "  +4.63%  "

$ws.Range("E38").Value = "  -3.22%  "

$ws.Range("D39").Value = "15.02"
$ws.Range("E39").Value = "  -7.85%  "

$ws.Range("D40").Value = "3.84"
$ws.Range("E40").Value = "  -4.93%  "

$ws.Range("E41").Value = "  -3.23%  "

$ws.Range("D42").Value = "3.20"
$ws.Range("E42").Value = "  -8.15%  "

$ws.Range("D43").Value = "1.939.34"
$ws.Range("E43").Value = "  +4.43%  "

$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").Value = "94.24"
$ws.Range("E45").Value = "  +6.12%  "

$ws.Range("D46").Value = "1.79"
$ws.Range("E46").Value = "  -11.63%  "

$ws.Range("D47").Value = "8.46"
$ws.Range("E47").Value = "  +5.48%  "

$ws.Range("E48").Value = "  -5.87%  "

$ws.Range("D49").Value = "98.90"
$ws.Range("E49").Value = "  +1.51%  "

$ws.Range("D50").Value = "2.602.86"
$ws.Range("E50").Value = "  +2.60%  "

$ws.Range("D51").Value = "68.51"
$ws.Range("E51").Value = "  -8.55%  "

# Restore the workbook default style now that the text is safely stored,
# so no stray number-format override is left on these cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
